# Laba_2: unique_lines.xlsx edit
# 1) Drop the three near-duplicate low-count rows that mixed the
#    "Косметика и медицина" category with the "Одежда,обувь и аксессуары"
#    brand (originally rows 12, 14 and 16 - the "от 5000 до 50000" bucket
#    for that odd combination). Delete from the bottom up so row indices
#    of the earlier deletions stay valid.
# 2) Rewrite the "Координаты и время" column (B) text for every
#    remaining data row to the corrected "lat, lon" formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(12).Delete()

$ws.Range("B2:B21").Value = "59, 30"
$ws.Range("B22:B29").Value = "59, 29"
$ws.Range("B30:B37").Value = "59, 30"
$ws.Range("B38:B45").Value = "60, 30"
